$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47, shifting rows 47:147 down to 48:148
$ws.Rows("47:47").Insert()

# Populate the newly inserted row 47 with the new data
$ws.Range("A47").Value = 10
$ws.Range("B47").Value = "Vega Modelo de Temuco"
$ws.Range("C47").Value = "La Araucanía"
$ws.Range("D47").Value = 44498
$ws.Range("E47").Value = 9
$ws.Range("F47").Value = 100112005
$ws.Range("G47").Value = "Puerro"
$ws.Range("H47").Value = "Azul de Maquehue"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 60
$ws.Range("K47").Value = 7000
$ws.Range("L47").Value = 7000
$ws.Range("M47").Value = 7000
$ws.Range("N47").Value = "$/docena de paquetes"
$ws.Range("O47").Value = "Provincia de Cautín"
$ws.Range("P47").Value = 583
$ws.Range("Q47").Value = 12
$ws.Range("R47").Value = "Hortaliza"
